# Weekly update: a new "Arveja Verde" price record is inserted as the
# 2nd data row (worksheet row 4), shifting all the following rows down
# by one. The new record's date, volume, prices, origin region and
# $/Kg values are populated; the remaining (constant-for-this-sheet)
# columns keep the values that were already present in that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 4 (old row 4 and everything
# below shifts down by one row).
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new record's data.
$ws.Range("A4").Value = 10
$ws.Range("B4").Value = "Vega Modelo de Temuco"
$ws.Range("C4").Value = "La Araucanía"
$ws.Range("D4").Value = 44515
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = 100112022
$ws.Range("G4").Value = "Arveja Verde"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 115
$ws.Range("K4").Value = 16000
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 16000
$ws.Range("N4").Value = "$/saco 25 kilos"
$ws.Range("O4").Value = "Región del Maule"
$ws.Range("P4").Value = 640
$ws.Range("Q4").Value = 25
$ws.Range("R4").Value = "Hortaliza"
